$d = $word.ActiveDocument

$targets = @("Alarmed by", "Recognizing", "Noting with concern")

foreach ($t in $targets) {
    $rng = $d.Content
    $found = $rng.Find.Execute($t, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Font.Italic = 1
        $rng.Font.Underline = 0
    }
}
